$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A13").Value = "Feat: Dev prod code to refresh games daily"
$ws.Range("B13").Value2 = $ws.Range("B2").Value2
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C13").Value = "M"

$ws.Range("A14").Select()
